$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 22:52"

# --- Update per-country statistics with freshly reported numbers ---
# Columns: A Pais | B Casos totales | C Nuevos casos | D Casos activos
#          E Recuperados | F Casos criticos | G Muertes hoy | H Muertes

# Estados Unidos (row 4)
$ws.Range("B4").Value = 876156
$ws.Range("C4").Value = 27439
$ws.Range("E4").Value = 741498
$ws.Range("F4").Value = 14994
$ws.Range("G4").Value = 1989
$ws.Range("H4").Value = 49648

# Brasil (row 14)
$ws.Range("B14").Value = 49492
$ws.Range("C14").Value = 3735
$ws.Range("E14").Value = 20861
$ws.Range("G14").Value = 407
$ws.Range("H14").Value = 3313

# Barein (row 65)
$ws.Range("B65").Value = 2217
$ws.Range("C65").Value = 190
$ws.Range("E65").Value = 1127

# Costa de Marfil (row 87 before re-sort)
$ws.Range("B87").Value = 1004
$ws.Range("C87").Value = 52
$ws.Range("D87").Value = 359
$ws.Range("E87").Value = 631
$ws.Range("H87").Value = 14

# Trinidad yTobago (row 141)
$ws.Range("D141").Value = 45
$ws.Range("E141").Value = 62

# Republica de Africa Central (row 192 before re-sort)
$ws.Range("B192").Value = 16
$ws.Range("D192").Value = 10
$ws.Range("E192").Value = 6

# --- Re-sort the country table by "Casos totales" (column B) descending,
#     exactly like the published sheet does every refresh, so the two
#     countries whose totals just overtook their neighbour (Costa de Marfil
#     over Republica de Yibuti, Republica de Africa Central over San
#     Cristobal y Nieves) move up into their new rank. ---
$dataRange = $ws.Range("A4:H216")
$key1 = $ws.Range("B4:B216")
$dataRange.Sort($key1, 2, $null, $null, 1, $null, 1, 1)
